$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data row (row 4) with the new projection figures
$ws.Range("A4").Value = 41610
$ws.Range("B4").Value = 37717
$ws.Range("C4").Value = 52147
$ws.Range("D4").Value = 1.849333333333334
$ws.Range("E4").Value = 1.676333333333334
$ws.Range("F4").Value = 2.317666666666667
$ws.Range("G4").Value = 0.3016666666666667
$ws.Range("H4").Value = 0.2726666666666667
$ws.Range("I4").Value = 0.3183333333333334
$ws.Range("J4").Value = 0.7914315642133092
$ws.Range("K4").Value = 0.7317640931035322
$ws.Range("L4").Value = 0.7426540931035323
$ws.Range("M4").Value = 435
$ws.Range("N4").Value = 292
$ws.Range("O4").Value = 1192
$ws.Range("P4").Value = 0.02023255813953489
$ws.Range("Q4").Value = 0.0136046511627907
$ws.Range("R4").Value = 0.05546511627906978
$ws.Range("S4").Value = 0.004883720930232559
$ws.Range("T4").Value = 0.004186046511627908
$ws.Range("U4").Value = 0.009069767441860466
$ws.Range("V4").Value = 0.7918915642133091
$ws.Range("W4").Value = 0.7320940931035322
$ws.Range("X4").Value = 0.7439790931035322
